# chore: adapt column header formatting to respective input file names
#
# Renames the header row suffixes ("_old" -> "_FV2210", "_new" -> "_FV2304"),
# freezes the header row, and turns the data range into a native Excel
# Table (ListObject) with an AutoFilter, matching the new "<formatversion>"
# header naming scheme.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header row (row 1, columns A:U) ----------------------------
$headerMap = @{
    "A1" = "Segmentname_FV2210";
    "B1" = "Segmentgruppe_FV2210";
    "C1" = "Segment_FV2210";
    "D1" = "Datenelement_FV2210";
    "E1" = "Segment ID_FV2210";
    "F1" = "Code_FV2210";
    "G1" = "Qualifier_FV2210";
    "H1" = "Beschreibung_FV2210";
    "I1" = "Bedingungsausdruck_FV2210";
    "J1" = "Bedingung_FV2210";
    "K1" = "diff";
    "L1" = "Segmentname_FV2304";
    "M1" = "Segmentgruppe_FV2304";
    "N1" = "Segment_FV2304";
    "O1" = "Datenelement_FV2304";
    "P1" = "Segment ID_FV2304";
    "Q1" = "Code_FV2304";
    "R1" = "Qualifier_FV2304";
    "S1" = "Beschreibung_FV2304";
    "T1" = "Bedingungsausdruck_FV2304";
    "U1" = "Bedingung_FV2304";
}

foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# --- 2) Freeze the header row (top row) -----------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3) Turn the used range into a Table (ListObject) with AutoFilter -----
$usedRange = $ws.Range("A1:U85")
$table = $ws.ListObjects.Add(1, $usedRange, $false, 1, "")
$table.Name = "Table1"
